function Set-CellText {
    param($ws, $addr, $text)
    # Force the cell to stay text (avoid Excel auto-converting numeric-looking
    # strings like "204.74" into a floating point number), then restore the
    # cell's original (default) style so no stray formatting is introduced.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-CellText $ws "D2" "79.075.12"
Set-CellText $ws "E2" "  +3.20%  "

# Row 3 - Ethereum
Set-CellText $ws "D3" "3.181.43"
Set-CellText $ws "E3" "  +5.05%  "

# Row 4 - TetherUSD
Set-CellText $ws "D4" "0.999"
Set-CellText $ws "E4" "  -0.10%  "

# Row 5 - Solana
Set-CellText $ws "D5" "204.74"
Set-CellText $ws "E5" "  +1.51%  "

# Row 6 - BNB
Set-CellText $ws "D6" "634.63"
Set-CellText $ws "E6" "  +0.56%  "

# Row 7 - USDC
Set-CellText $ws "D7" "0.999"
Set-CellText $ws "E7" "  -0.03%  "

# Row 8 - Dogecoin
Set-CellText $ws "D8" "0.230"
Set-CellText $ws "E8" "  +9.20%  "

# Row 9 - XRP
Set-CellText $ws "D9" "0.583"
Set-CellText $ws "E9" "  +5.09%  "

# Row 10 - LidoStakedEther
Set-CellText $ws "D10" "3.176.62"
Set-CellText $ws "E10" "  +4.90%  "

# Row 11 - Cardano
Set-CellText $ws "D11" "0.580"
Set-CellText $ws "E11" "  +32.78%  "

# Row 12 - TRON
Set-CellText $ws "E12" "  +2.97%  "

# Row 13 - Toncoin
Set-CellText $ws "E13" "  +4.19%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-CellText $ws "D14" "3.765.17"
Set-CellText $ws "E14" "  +4.99%  "

# Row 15 - ShibaInu
Set-CellText $ws "D15" "0.0000226"
Set-CellText $ws "E15" "  +16.28%  "

# Row 16 - Avalanche
Set-CellText $ws "D16" "31.49"
Set-CellText $ws "E16" "  +6.96%  "

# Row 17 - WrappedBTC
Set-CellText $ws "D17" "78.874.31"
Set-CellText $ws "E17" "  +3.04%  "

# Row 18 - WrappedEther
Set-CellText $ws "D18" "3.178.15"
Set-CellText $ws "E18" "  +5.07%  "

# Row 19 - Chainlink
Set-CellText $ws "D19" "14.44"
Set-CellText $ws "E19" "  +7.07%  "

# Rows 20/21 swap: Uniswap <-> SuiNetwork
Set-CellText $ws "B20" "SuiNetwork"
Set-CellText $ws "C20" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText $ws "D20" "2.98"
Set-CellText $ws "E20" "  +29.19%  "

Set-CellText $ws "B21" "Uniswap"
Set-CellText $ws "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText $ws "D21" "9.25"
Set-CellText $ws "E21" "  +2.23%  "

# Row 22 - BitcoinCash
Set-CellText $ws "D22" "426.60"
Set-CellText $ws "E22" "  +12.97%  "

# Row 23 - Polkadot
Set-CellText $ws "D23" "4.98"
Set-CellText $ws "E23" "  +13.69%  "

# Row 24 - LEO
Set-CellText $ws "D24" "6.86"
Set-CellText $ws "E24" "  +5.79%  "

# Row 25 - NEARProtocol
Set-CellText $ws "D25" "4.78"
Set-CellText $ws "E25" "  +8.11%  "

# Row 26 - Aptos
Set-CellText $ws "D26" "11.19"
Set-CellText $ws "E26" "  +12.57%  "

# Row 27 - Litecoin
Set-CellText $ws "D27" "76.44"
Set-CellText $ws "E27" "  +3.72%  "

# Row 28 - Dai
Set-CellText $ws "E28" "  +0.07%  "

# Row 29 - PEPE
Set-CellText $ws "D29" "0.0000115"
Set-CellText $ws "E29" "  +3.33%  "

# Rows 30/31 swap: Binance-PegBSC-USD <-> InternetComputer(DFINITY)
Set-CellText $ws "B30" "InternetComputer(DFINITY)"
Set-CellText $ws "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText $ws "D30" "8.98"
Set-CellText $ws "E30" "  +7.54%  "

Set-CellText $ws "B31" "Binance-PegBSC-USD"
Set-CellText $ws "C31" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-CellText $ws "D31" "0.992"
Set-CellText $ws "E31" "  -0.51%  "

# Row 32 - Fetch.AI
Set-CellText $ws "D32" "1.48"
Set-CellText $ws "E32" "  +4.27%  "

# Row 33 - Bittensor
Set-CellText $ws "D33" "520.75"
Set-CellText $ws "E33" "  +1.77%  "

# Row 34 - PancakeSwap
Set-CellText $ws "D34" "2.01"
Set-CellText $ws "E34" "  +2.23%  "

# Row 35 - Kaspa
Set-CellText $ws "D35" "0.138"
Set-CellText $ws "E35" "  +21.77%  "

# Row 36 - EthereumClassic
Set-CellText $ws "D36" "22.96"
Set-CellText $ws "E36" "  +10.58%  "

# Row 37 - Cronos
Set-CellText $ws "E37" "  +15.32%  "

# Row 38 - FirstDigitalUSD
Set-CellText $ws "D38" "0.999"
Set-CellText $ws "E38" "  -0.13%  "

# Row 39 - PolygonEcosystemToken
Set-CellText $ws "D39" "0.402"
Set-CellText $ws "E39" "  +4.60%  "

# Row 40 - Monero
Set-CellText $ws "D40" "163.97"
Set-CellText $ws "E40" "  +0.48%  "

# Row 41 - WhiteBITCoin
Set-CellText $ws "D41" "19.99"
Set-CellText $ws "E41" "  -0.18%  "

# Row 42 - USDe
Set-CellText $ws "E42" "  +0.03%  "

# Row 43 - Aave
Set-CellText $ws "D43" "191.74"
Set-CellText $ws "E43" "  +1.75%  "

# Row 44 - RenderToken
Set-CellText $ws "D44" "5.44"
Set-CellText $ws "E44" "  +4.98%  "

# Row 45 - Mantle
Set-CellText $ws "D45" "0.821"
Set-CellText $ws "E45" "  +14.06%  "

# Row 46 - Stacks
Set-CellText $ws "D46" "1.79"
Set-CellText $ws "E46" "  +6.69%  "

# Row 47 - ImmutableX
Set-CellText $ws "E47" "  +3.42%  "

# Row 48 - OKB
Set-CellText $ws "D48" "42.54"
Set-CellText $ws "E48" "  +0.14%  "

# Rows 49/50 swap: InjectiveProtocol <-> dogwifhat
Set-CellText $ws "B49" "dogwifhat"
Set-CellText $ws "C49" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-CellText $ws "D49" "2.51"
Set-CellText $ws "E49" "  +2.13%  "

Set-CellText $ws "B50" "InjectiveProtocol"
Set-CellText $ws "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText $ws "D50" "25.28"
Set-CellText $ws "E50" "  +11.84%  "

# Row 51 - ARBITRUM
Set-CellText $ws "D51" "0.624"
Set-CellText $ws "E51" "  +2.64%  "
